$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revised (reseasonalized) values for existing rows 90-142, columns B:G ---
# Flat list, 6 numbers per original row (B,C,D,E,F,G), row 90 first.
$flatVals = @(
    26507,6357,20125,10819,8766,28609,
    27536,6682,20820,10279,9063,28626,
    27063,6666,20353,10358,8955,28386,
    26264,6397,19831,10011,8130,28115,
    24822,5915,18886,9920,7131,27879,
    24594,5588,19012,9714,6896,27710,
    24880,5490,19411,9919,7067,28021,
    26109,5650,20492,10163,8000,28376,
    26924,5991,20953,9824,8418,28471,
    28313,6236,22105,9872,8869,29443,
    28880,6425,22476,10709,9548,30132,
    29877,6871,23002,10244,9693,30510,
    30632,7316,23318,10579,10394,30891,
    31135,7251,23892,10712,10595,31321,
    31546,7296,24259,10650,10810,31435,
    31378,7814,23553,10903,10323,32045,
    32254,7700,24556,10788,10497,32597,
    33200,7763,25441,10715,10836,33089,
    34350,8582,25768,10198,11167,33294,
    34081,8883,25195,11221,11783,33533,
    34547,8666,25881,10619,11242,33924,
    34623,8808,25816,11387,11563,34448,
    34684,8609,26075,11406,11347,34742,
    34869,8191,26678,11012,11065,34815,
    34538,8262,26275,11075,10746,34867,
    34283,8331,25952,11324,10576,35031,
    34295,7930,26364,10938,10230,35002,
    34907,8044,26860,11214,10715,35406,
    34872,8039,26836,11038,10418,35500,
    35144,7981,27170,10784,9991,35965,
    35695,8349,27347,10852,10704,35830,
    35756,8080,27684,11070,10645,36180,
    35802,8105,27707,11099,10489,36418,
    35644,8139,27511,10936,10252,36341,
    36293,7992,28320,10916,10751,36473,
    36211,7782,28457,10992,10647,36569,
    36765,7879,28914,10446,11020,36268,
    36957,7605,29391,10511,10889,36646,
    37021,7708,29349,11138,10881,37303,
    37718,7843,29912,11163,11294,37635,
    38092,8019,30107,11489,11594,38055,
    38931,8073,30896,11445,12089,38386,
    38684,8130,30589,11285,11778,38272,
    39398,8379,31052,11395,12189,38706,
    39094,8379,30747,11250,11831,38599,
    39652,8486,31199,11008,11724,39019,
    39595,8566,31059,11314,11643,39331,
    38186,8625,29582,10886,11325,37818,
    37878,8320,29582,11309,10671,38552,
    31738,6791,24977,11259,9334,33665,
    34784,7095,27742,10580,9904,35493,
    37595,7875,29766,10811,10669,37788,
    40931,8461,32525,10647,12508,39460
)

$nRows = 53
$nCols = 6
$updates = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $updates[$i,$j] = $flatVals[$i * $nCols + $j]
    }
}
$ws.Range("B90:G142").Value = $updates

# --- New row 143: quarter 01-04-2021 ---
$newRow = New-Object 'object[,]' 1,6
$newRow[0,0] = 41571
$newRow[0,1] = 8487
$newRow[0,2] = 33147
$newRow[0,3] = 10919
$newRow[0,4] = 13043
$newRow[0,5] = 39460
$ws.Range("B143:G143").Value = $newRow

# Column A label "01-04-2021" must stay literal text (not auto-parsed into a date
# serial the way a plain `.Value = "01-04-2021"` assignment would be). Routing it
# through a formula result first keeps it a string, then Copy + PasteSpecial
# (values-only) bakes it in as a plain shared-string cell, matching the style-free
# "t=s" cells used by every other date label in column A.
$aCell = $ws.Cells.Item(143, 1)
$aCell.Formula = '="01-04-2021"'
$aCell.Copy()
$aCell.PasteSpecial(-4163)
$excel.CutCopyMode = $false
